# Add columns I (I0) and J (IF) to the worksheet, matching the header style
# used by the existing columns (A..H), and populate all 58 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell formatting (bold font, borders, center/top alignment)
# from the existing "IP" header (H1) onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# row, I-value, J-value
$data = @(
    @(2, 9, 9),
    @(3, 10, 11),
    @(4, 4, 5),
    @(5, 9, 9),
    @(6, 8, 8),
    @(7, 9, 9),
    @(8, 5, 5),
    @(9, 9, 9),
    @(10, 9, 9),
    @(11, 6, 7),
    @(12, 8, 8),
    @(13, 8, 8),
    @(14, 6, 6),
    @(15, 7, 8),
    @(16, 9, 9),
    @(17, 4, 5),
    @(18, 6, 7),
    @(19, 7, 7),
    @(20, 9, 9),
    @(21, 11, 11),
    @(22, 7, 8),
    @(23, 8, 9),
    @(24, 7, 7),
    @(25, 7, 7),
    @(26, 9, 9),
    @(27, 8, 8),
    @(28, 12, 13),
    @(29, 8, 8),
    @(30, 8, 9),
    @(31, 6, 8),
    @(32, 3, 5),
    @(33, 8, 8),
    @(34, 5, 6),
    @(35, 8, 8),
    @(36, 5, 6),
    @(37, 7, 7),
    @(38, 7, 8),
    @(39, 6, 7),
    @(40, 7, 7),
    @(41, 9, 9),
    @(42, 7, 8),
    @(43, 4, 5),
    @(44, 9, 9),
    @(45, 6, 6),
    @(46, 6, 6),
    @(47, 6, 7),
    @(48, 9, 9),
    @(49, 8, 8),
    @(50, 7, 7),
    @(51, 5, 6),
    @(52, 8, 8),
    @(53, 9, 9),
    @(54, 9, 9),
    @(55, 5, 6),
    @(56, 8, 8),
    @(57, 9, 9),
    @(58, 4, 5),
    @(59, 4, 4),
)

foreach ($entry in $data) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
